$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header summary values
$ws.Range("E11").Value = 455520
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 8

# Update detail rows 16-22: periods now run in ascending order 2501..2507
$ws.Range("E16").Value = "2501"
$ws.Range("E17").Value = "2502"
$ws.Range("E18").Value = "2503"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2505"
$ws.Range("E21").Value = "2506"
$ws.Range("E22").Value = "2507"

# Row 23: now belongs to the same worker (BETTY PALOMINO ALTAMIRANDA), period 2508
$ws.Range("C23").Value = "1047375786"
$ws.Range("D23").Value = "BETTY PALOMINO ALTAMIRANDA"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
